$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "304.06"
Set-TextValue $ws.Range("E2") "3.75%"
Set-TextValue $ws.Range("G2") "23"
Set-TextValue $ws.Range("D3") "32.06"
Set-TextValue $ws.Range("E3") "6.87%"
Set-TextValue $ws.Range("G3") "23"
Set-TextValue $ws.Range("D4") "5.253"
Set-TextValue $ws.Range("E4") "2.17%"
Set-TextValue $ws.Range("G4") "23"
Set-TextValue $ws.Range("D5") "0.07559"
Set-TextValue $ws.Range("E5") "5.87%"
Set-TextValue $ws.Range("G5") "23"
Set-TextValue $ws.Range("D6") "7.873"
Set-TextValue $ws.Range("E6") "4.27%"
Set-TextValue $ws.Range("G6") "23"
Set-TextValue $ws.Range("D7") "3.859"
Set-TextValue $ws.Range("E7") "6.51%"
Set-TextValue $ws.Range("G7") "23"
Set-TextValue $ws.Range("D8") "1.559"
Set-TextValue $ws.Range("E8") "11.32%"
Set-TextValue $ws.Range("G8") "23"
Set-TextValue $ws.Range("D9") "0.9298"
Set-TextValue $ws.Range("E9") "1.51%"
Set-TextValue $ws.Range("G9") "23"
Set-TextValue $ws.Range("D10") "0.1688"
Set-TextValue $ws.Range("E10") "3.51%"
Set-TextValue $ws.Range("G10") "23"
Set-TextValue $ws.Range("D11") "0.07979"
Set-TextValue $ws.Range("E11") "4.16%"
Set-TextValue $ws.Range("G11") "23"
Set-TextValue $ws.Range("D12") "0.08018"
Set-TextValue $ws.Range("E12") "3.49%"
Set-TextValue $ws.Range("G12") "23"
Set-TextValue $ws.Range("D13") "0.03033"
Set-TextValue $ws.Range("E13") "3.38%"
Set-TextValue $ws.Range("G13") "23"
Set-TextValue $ws.Range("D14") "0.09921"
Set-TextValue $ws.Range("E14") "10.20%"
Set-TextValue $ws.Range("G14") "23"
Set-TextValue $ws.Range("D15") "0.001488"
Set-TextValue $ws.Range("E15") "-7.41%"
Set-TextValue $ws.Range("G15") "23"
Set-TextValue $ws.Range("D16") "0.04599"
Set-TextValue $ws.Range("E16") "1.28%"
Set-TextValue $ws.Range("G16") "23"
Set-TextValue $ws.Range("D17") "0.006229"
Set-TextValue $ws.Range("E17") "-1.88%"
Set-TextValue $ws.Range("G17") "23"
Set-TextValue $ws.Range("D18") "3.446"
Set-TextValue $ws.Range("E18") "-1.10%"
Set-TextValue $ws.Range("G18") "23"
Set-TextValue $ws.Range("D19") "2.238"
Set-TextValue $ws.Range("E19") "0.19%"
Set-TextValue $ws.Range("G19") "23"
Set-TextValue $ws.Range("E20") "0.88%"
Set-TextValue $ws.Range("G20") "23"
Set-TextValue $ws.Range("D21") "0.1333"
Set-TextValue $ws.Range("E21") "-2.61%"
Set-TextValue $ws.Range("G21") "23"
Set-TextValue $ws.Range("D22") "4.549"
Set-TextValue $ws.Range("E22") "13.10%"
Set-TextValue $ws.Range("G22") "23"
Set-TextValue $ws.Range("D23") "0.1614"
Set-TextValue $ws.Range("E23") "1.06%"
Set-TextValue $ws.Range("G23") "23"
Set-TextValue $ws.Range("D24") "0.001211"
Set-TextValue $ws.Range("E24") "0.17%"
Set-TextValue $ws.Range("G24") "23"
Set-TextValue $ws.Range("D25") "0.004482"
Set-TextValue $ws.Range("E25") "5.21%"
Set-TextValue $ws.Range("G25") "23"
Set-TextValue $ws.Range("D26") "0.0001393"
Set-TextValue $ws.Range("E26") "19.16%"
Set-TextValue $ws.Range("G26") "23"
Set-TextValue $ws.Range("D27") "0.0001778"
Set-TextValue $ws.Range("E27") "7.66%"
Set-TextValue $ws.Range("G27") "23"
Set-TextValue $ws.Range("G28") "23"
Set-TextValue $ws.Range("G29") "23"
Set-TextValue $ws.Range("G30") "23"
Set-TextValue $ws.Range("G31") "23"
Set-TextValue $ws.Range("G32") "23"
Set-TextValue $ws.Range("G33") "23"
Set-TextValue $ws.Range("G34") "23"
Set-TextValue $ws.Range("G35") "23"
Set-TextValue $ws.Range("G36") "23"
Set-TextValue $ws.Range("G37") "23"
Set-TextValue $ws.Range("G38") "23"
Set-TextValue $ws.Range("D39") "0.01726"
Set-TextValue $ws.Range("E39") "2,537.97%"
Set-TextValue $ws.Range("G39") "23"
Set-TextValue $ws.Range("D40") "0.04494"
Set-TextValue $ws.Range("E40") "2.09%"
Set-TextValue $ws.Range("G40") "23"
Set-TextValue $ws.Range("D41") "0.006905"
Set-TextValue $ws.Range("E41") "-1.72%"
Set-TextValue $ws.Range("G41") "23"
Set-TextValue $ws.Range("E42") "6.31%"
Set-TextValue $ws.Range("G42") "23"
Set-TextValue $ws.Range("D43") "0.002070"
Set-TextValue $ws.Range("E43") "-6.28%"
Set-TextValue $ws.Range("G43") "23"
Set-TextValue $ws.Range("D44") "0.01370"
Set-TextValue $ws.Range("E44") "2.68%"
Set-TextValue $ws.Range("G44") "23"
Set-TextValue $ws.Range("D45") "0.00006143"
Set-TextValue $ws.Range("E45") "4.92%"
Set-TextValue $ws.Range("G45") "23"
Set-TextValue $ws.Range("G46") "23"
Set-TextValue $ws.Range("D47") "0.01294"
Set-TextValue $ws.Range("E47") "-0.22%"
Set-TextValue $ws.Range("G47") "23"
Set-TextValue $ws.Range("G48") "23"
Set-TextValue $ws.Range("G49") "23"
Set-TextValue $ws.Range("G50") "23"
Set-TextValue $ws.Range("G51") "23"
